$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ZZ1 is a scratch cell (outside the used A1:E51 range) used to stage
# text values for numeric-looking strings. Writing straight into a
# General-formatted destination cell would auto-convert a string like
# "246.27" into a Number; formatting ZZ1 as Text, writing the string
# there, then Copy + PasteSpecial(xlPasteValues=-4163) into the real
# destination preserves the value as TEXT (matching the source file,
# which stores these as inline strings) without leaving any extra
# number-format / style applied to the destination cell itself.

$ws.Range('D2').Value = '37.414.31'
$ws.Range('E2').Value = '  +2.41%  '
$ws.Range('D3').Value = '2.004.10'
$ws.Range('E3').Value = '  +2.21%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('ZZ1').NumberFormat = '@'
$ws.Range('ZZ1').Value = '246.27'
$ws.Range('ZZ1').Copy()
$ws.Range('D5').PasteSpecial(-4163)
$ws.Range('E5').Value = '  +0.70%  '
$ws.Range('E6').Value = '  +2.57%  '
$ws.Range('ZZ1').NumberFormat = '@'
$ws.Range('ZZ1').Value = '61.89'
$ws.Range('ZZ1').Copy()
$ws.Range('D7').PasteSpecial(-4163)
$ws.Range('E7').Value = '  +5.54%  '
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('ZZ1').NumberFormat = '@'
$ws.Range('ZZ1').Value = '0.0804'
$ws.Range('ZZ1').Copy()
$ws.Range('D10').PasteSpecial(-4163)
$ws.Range('E10').Value = '  -0.32%  '
$ws.Range('E11').Value = '  +0.37%  '
$ws.Range('ZZ1').NumberFormat = '@'
$ws.Range('ZZ1').Value = '14.94'
$ws.Range('ZZ1').Copy()
$ws.Range('D12').PasteSpecial(-4163)
$ws.Range('E12').Value = '  +8.69%  '
$ws.Range('ZZ1').NumberFormat = '@'
$ws.Range('ZZ1').Value = '22.63'
$ws.Range('ZZ1').Copy()
$ws.Range('D13').PasteSpecial(-4163)
$ws.Range('E13').Value = '  +1.97%  '
$ws.Range('ZZ1').NumberFormat = '@'
$ws.Range('ZZ1').Value = '0.852'
$ws.Range('ZZ1').Copy()
$ws.Range('D14').PasteSpecial(-4163)
$ws.Range('E14').Value = '  +2.52%  '
$ws.Range('D15').Value = '2.298.42'
$ws.Range('E15').Value = '  +2.21%  '
$ws.Range('ZZ1').NumberFormat = '@'
$ws.Range('ZZ1').Value = '5.45'
$ws.Range('ZZ1').Copy()
$ws.Range('D16').PasteSpecial(-4163)
$ws.Range('E16').Value = '  +2.85%  '
$ws.Range('D17').Value = '2.009.18'
$ws.Range('E17').Value = '  +2.67%  '
$ws.Range('D18').Value = '37.331.64'
$ws.Range('E18').Value = '  +2.24%  '
$ws.Range('ZZ1').NumberFormat = '@'
$ws.Range('ZZ1').Value = '70.53'
$ws.Range('ZZ1').Copy()
$ws.Range('D19').PasteSpecial(-4163)
$ws.Range('E19').Value = '  +1.07%  '
$ws.Range('D20').Value = '0.0₃0868'
$ws.Range('E20').Value = '  +1.47%  '
$ws.Range('E21').Value = '  +3.52%  '
$ws.Range('E23').Value = '  +0.00%  '
$ws.Range('ZZ1').NumberFormat = '@'
$ws.Range('ZZ1').Value = '2.52'
$ws.Range('ZZ1').Copy()
$ws.Range('D24').PasteSpecial(-4163)
$ws.Range('E24').Value = '  +2.88%  '
$ws.Range('ZZ1').NumberFormat = '@'
$ws.Range('ZZ1').Value = '2.36'
$ws.Range('ZZ1').Copy()
$ws.Range('D25').PasteSpecial(-4163)
$ws.Range('E25').Value = '  +0.40%  '
$ws.Range('ZZ1').NumberFormat = '@'
$ws.Range('ZZ1').Value = '0.145'
$ws.Range('ZZ1').Copy()
$ws.Range('D26').PasteSpecial(-4163)
$ws.Range('E26').Value = '  +5.17%  '
$ws.Range('ZZ1').NumberFormat = '@'
$ws.Range('ZZ1').Value = '9.35'
$ws.Range('ZZ1').Copy()
$ws.Range('D27').PasteSpecial(-4163)
$ws.Range('E27').Value = '  +1.12%  '
$ws.Range('ZZ1').NumberFormat = '@'
$ws.Range('ZZ1').Value = '163.78'
$ws.Range('ZZ1').Copy()
$ws.Range('D28').PasteSpecial(-4163)
$ws.Range('E28').Value = '  +2.05%  '
$ws.Range('ZZ1').NumberFormat = '@'
$ws.Range('ZZ1').Value = '19.76'
$ws.Range('ZZ1').Copy()
$ws.Range('D29').PasteSpecial(-4163)
$ws.Range('E29').Value = '  +1.68%  '
$ws.Range('ZZ1').NumberFormat = '@'
$ws.Range('ZZ1').Value = '1.37'
$ws.Range('ZZ1').Copy()
$ws.Range('D30').PasteSpecial(-4163)
$ws.Range('E30').Value = '  +19.17%  '
$ws.Range('E31').Value = '  +1.54%  '
$ws.Range('ZZ1').NumberFormat = '@'
$ws.Range('ZZ1').Value = '4.90'
$ws.Range('ZZ1').Copy()
$ws.Range('D32').PasteSpecial(-4163)
$ws.Range('E32').Value = '  +4.21%  '
$ws.Range('ZZ1').NumberFormat = '@'
$ws.Range('ZZ1').Value = '0.0626'
$ws.Range('ZZ1').Copy()
$ws.Range('D33').PasteSpecial(-4163)
$ws.Range('E33').Value = '  +0.96%  '
$ws.Range('ZZ1').NumberFormat = '@'
$ws.Range('ZZ1').Value = '4.61'
$ws.Range('ZZ1').Copy()
$ws.Range('D34').PasteSpecial(-4163)
$ws.Range('E34').Value = '  +6.74%  '
$ws.Range('ZZ1').NumberFormat = '@'
$ws.Range('ZZ1').Value = '2.35'
$ws.Range('ZZ1').Copy()
$ws.Range('D35').PasteSpecial(-4163)
$ws.Range('E35').Value = '  +4.76%  '
$ws.Range('E36').Value = '  -0.06%  '
$ws.Range('E37').Value = '  +2.12%  '
$ws.Range('ZZ1').NumberFormat = '@'
$ws.Range('ZZ1').Value = '3.36'
$ws.Range('ZZ1').Copy()
$ws.Range('D38').PasteSpecial(-4163)
$ws.Range('E38').Value = '  -1.01%  '
$ws.Range('ZZ1').NumberFormat = '@'
$ws.Range('ZZ1').Value = '5.53'
$ws.Range('ZZ1').Copy()
$ws.Range('D39').PasteSpecial(-4163)
$ws.Range('E39').Value = '  -3.44%  '
$ws.Range('ZZ1').NumberFormat = '@'
$ws.Range('ZZ1').Value = '0.0983'
$ws.Range('ZZ1').Copy()
$ws.Range('D40').PasteSpecial(-4163)
$ws.Range('E40').Value = '  +0.19%  '
$ws.Range('E41').Value = '  +1.44%  '
$ws.Range('E42').Value = '  +1.95%  '
$ws.Range('E43').Value = '  +1.26%  '
$ws.Range('ZZ1').NumberFormat = '@'
$ws.Range('ZZ1').Value = '16.76'
$ws.Range('ZZ1').Copy()
$ws.Range('D44').PasteSpecial(-4163)
$ws.Range('E44').Value = '  +4.58%  '
$ws.Range('ZZ1').NumberFormat = '@'
$ws.Range('ZZ1').Value = '91.37'
$ws.Range('ZZ1').Copy()
$ws.Range('D45').PasteSpecial(-4163)
$ws.Range('E45').Value = '  +3.94%  '
$ws.Range('D46').Value = '1.385.48'
$ws.Range('E46').Value = '  +1.43%  '
$ws.Range('E47').Value = '  +0.98%  '
$ws.Range('E48').Value = '  +1.65%  '
$ws.Range('B49').Value = 'MultiversX'
$ws.Range('C49').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('ZZ1').NumberFormat = '@'
$ws.Range('ZZ1').Value = '46.78'
$ws.Range('ZZ1').Copy()
$ws.Range('D49').PasteSpecial(-4163)
$ws.Range('E49').Value = '  +6.87%  '
$ws.Range('B50').Value = 'MXToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('ZZ1').NumberFormat = '@'
$ws.Range('ZZ1').Value = '2.84'
$ws.Range('ZZ1').Copy()
$ws.Range('D50').PasteSpecial(-4163)
$ws.Range('E50').Value = '  +0.41%  '
$ws.Range('E51').Value = '  +12.98%  '

# Clean up the scratch cell so it doesn't linger in the saved workbook.
$ws.Range('ZZ1').Clear()
